# Weekly price update: insert one new observation row for
# "Vega Modelo de Temuco - Cebollín" ahead of the existing row 219,
# pushing the existing data (rows 219-254) down by one row
# (new last row becomes 255), and fill the newly inserted row
# with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 219; Excel shifts rows 219:254 down to 220:255
# and carries the column D date-number-format down onto the new row.
$ws.Rows.Item(219).Insert()

$ws.Cells.Item(219, 1).Value  = 10
$ws.Cells.Item(219, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(219, 3).Value  = "La Araucanía"
$ws.Cells.Item(219, 4).Value  = 44522
$ws.Cells.Item(219, 5).Value  = 9
$ws.Cells.Item(219, 6).Value  = 100112037
$ws.Cells.Item(219, 7).Value  = "Cebollín"
$ws.Cells.Item(219, 8).Value  = "Sin especificar"
$ws.Cells.Item(219, 9).Value  = "Primera"
$ws.Cells.Item(219, 10).Value = 70
$ws.Cells.Item(219, 11).Value = 8000
$ws.Cells.Item(219, 12).Value = 9000
$ws.Cells.Item(219, 13).Value = 8571
$ws.Cells.Item(219, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(219, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(219, 16).Value = 714
$ws.Cells.Item(219, 17).Value = 12
$ws.Cells.Item(219, 18).Value = "Hortaliza"
